# Add an optional BLE121 "Case Print/Machining" line item to the BOM table.
# This is modeled as inserting a new row right before the existing
# "LightPipe" row (row 21), which pushes LightPipe and everything below it
# down by one row, then filling in the new row's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 21 (shifts LightPipe/blank/PCB/Assembly/Total
# rows down by one).
$ws.Rows("21:21").Insert()

# Grow the table (ListObject) to cover the newly inserted row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:J26"))

# Populate the new row with the BLE121 case print/machining line item.
$ws.Range("A21").Value = "Case Print/Machining"
$ws.Range("B21").Value = "PolyCase"
$ws.Range("F21").Value = 1
$ws.Range("I21").Value = 4.8899999999999997
$ws.Range("J21").Formula = "=F21*I21"

# The PCB row (now row 24) loses its subtotal formula -- it becomes a blank
# placeholder row just like the blank row above it, while the Assembly row
# (now row 25) keeps its own subtotal formula.
$ws.Range("J24").ClearContents()

# Match the author's final cursor position.
[void]$ws.Range("F22").Select()
